$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

# Update the four distribution fit strings for Week 16 (row 2: R, P, RA, PA)
$ws.Range("B2").Value = "JSU(-1.002473958785314, 1.0376526008197926, 0.5664992829547559, 2.70894152035578)"
$ws.Range("C2").Value = "JSU(-1.4154384521439676, 1.1889270007873973, 2.4168454981494905, 4.233756044135445)"
$ws.Range("D2").Value = "JSU(-0.9639947332828787, 1.1224219643730122, 0.6037600124686713, 2.6108148654236505)"
$ws.Range("E2").Value = "NCT(3.2318354626547663, 1.6426684988699376, -0.020866179492713, 4.7661895453193)"
